$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Values for column I (I0) and J (IF) for rows 2-30
$values = @(
    @(8,8),
    @(9,9),
    @(8,8),
    @(9,9),
    @(7,7),
    @(8,8),
    @(6,6),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(6,6),
    @(8,9),
    @(7,7),
    @(7,8),
    @(1,1),
    @(9,9),
    @(6,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(6,6),
    @(6,6)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
